# Matches.xlsx: the C3 "video" cell held a descriptive label ("Tiznados vs USMP")
# while its hyperlink silently pointed at the USMP match video and carried an
# explicit display-text override ("display=..." in the OOXML). This edit makes
# the cell text itself the video URL (matching how the C2 row already works),
# so the override is no longer needed, and it adjusts the saved UI selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Cell C3 ("video" column) now shows the YouTube link directly instead of
#    the old "Tiznados vs USMP" label.
$ws.Range("C3").Value = "https://www.youtube.com/watch?v=vmWdqhzpwTA"

# 2) The hyperlink on C3 already points at that same URL - just drop its
#    pinned TextToDisplay override so it once again follows the cell text
#    (same as the C2 hyperlink, which has no override either).
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = $null
    }
}

# 3) Move the active selection to H13 (matches the saved UI state).
$ws.Range("H13").Select() | Out-Null

Write-Output "done"
